$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45186 (2023-09-17) to 45188 (2023-09-19) for every data row
# (rows 2 through 504).
$ws.Range("C2:C504").Value = 45188
